$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.070.11"
$ws.Range("E2").Value = "  -0.60%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.893.64"
$ws.Range("E4").Value = "  +0.15%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "306.50"
$ws.Range("E5").Value = "  -0.36%  "
$ws.Range("E6").Value = "  +0.18%  "
$ws.Range("E7").Value = "  -0.49%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3758"
$ws.Range("E8").Value = "  -0.79%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07256"
$ws.Range("E9").Value = "  -0.31%  "
$ws.Range("E10").Value = "  -1.12%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.8984"
$ws.Range("E11").Value = "  -0.19%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08173"
$ws.Range("E12").Value = "  -0.88%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.937.40"
$ws.Range("E13").Value = "  +1.65%  "
$ws.Range("E14").Value = "  +1.03%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.294"
$ws.Range("E15").Value = "  +0.09%  "
$ws.Range("E16").Value = "  +0.13%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008572"
$ws.Range("E17").Value = "  -0.76%  "
$ws.Range("E18").Value = "  +0.43%  "
$ws.Range("E19").Value = "  +0.30%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "27.106.41"
$ws.Range("E20").Value = "  -0.70%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.082"
$ws.Range("E21").Value = "  +0.15%  "
$ws.Range("E22").Value = "  +0.59%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.403"
$ws.Range("E23").Value = "  -0.72%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "148.67"
$ws.Range("E24").Value = "  +1.82%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.290"
$ws.Range("E25").Value = "  -0.93%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "18.17"
$ws.Range("E26").Value = "  -0.06%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.736"
$ws.Range("E27").Value = "  -0.66%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "115.00"
$ws.Range("E28").Value = "  +0.05%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "4.781"
$ws.Range("E29").Value = "  -0.72%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.848"
$ws.Range("E30").Value = "  -3.03%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.09214"
$ws.Range("E31").Value = "  -0.21%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.05026"
$ws.Range("E32").Value = "  -0.59%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.7868"
$ws.Range("E33").Value = "  -2.27%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.211"
$ws.Range("E34").Value = "  -2.64%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.422"
$ws.Range("E35").Value = "  +2.96%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.967"
$ws.Range("E36").Value = "  -1.39%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.612"
$ws.Range("E37").Value = "  +1.00%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.5707"
$ws.Range("E38").Value = "  -0.39%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01989"
$ws.Range("E39").Value = "  +0.21%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.074"
$ws.Range("E40").Value = "  -0.43%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "9.009"
$ws.Range("E41").Value = "  +0.27%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.544"
$ws.Range("E42").Value = "  -1.46%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "115.90"
$ws.Range("E43").Value = "  -3.02%  "
$ws.Range("E44").Value = "  -0.23%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.4843"
$ws.Range("E45").Value = "  -0.28%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.002"
$ws.Range("E46").Value = "  +0.21%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "10.08"
$ws.Range("E47").Value = "  -1.76%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.619"
$ws.Range("E48").Value = "  -0.41%  "
$ws.Range("E49").Value = "  +1.36%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "63.43"
$ws.Range("E50").Value = "  -0.57%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05927"
$ws.Range("E51").Value = "  -0.30%  "
